$d = $word.ActiveDocument

# Locate the run boundary right before "proyecto desde la terminal de
# PowerShell..." (i.e. the end of the preceding "nuevo " run) so that we
# can keep it from being swallowed into the run we are about to edit.
$nuevoFind = $d.Content
[void]$nuevoFind.Find.Execute("nuevo ", $true, $false, $false, $false, $false, $true, 1, $false, "", 0)
$nuevoStart = $nuevoFind.Start
$nuevoEnd = $nuevoFind.End

# Locate "de realizar" - the fragment that becomes "se realizará".
$target = $d.Content
[void]$target.Find.Execute("de realizar", $true, $false, $false, $false, $false, $true, 1, $false, "", 0)
$start = $target.Start
$end = $target.End

# Replace "de realizar" with "se realizará". (This also happens to merge
# the preceding "nuevo " run into the edited run, and leaves the whole
# sentence as a single run - both get fixed below.)
$newText = "se realizará"
$editRange = $d.Range($start, $end)
$editRange.Text = $newText

# Force the new text to live in its own run, distinct from the text
# before and after it, by toggling (and restoring) a character format on
# exactly that span - this causes the engine to keep it split out instead
# of re-coalescing it with its neighbours.
$newRange = $d.Range($start, $start + $newText.Length)
$newRange.Bold = 1
$newRange.Bold = 0

# Likewise, re-isolate "nuevo " so it doesn't stay fused with the run
# that follows it.
$nuevoRange = $d.Range($nuevoStart, $nuevoEnd)
$nuevoRange.Bold = 1
$nuevoRange.Bold = 0
